$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 19
$ws1.Range("F5").Value = 20
$ws1.Range("F8").Value = 20
$ws1.Range("F9").Value = 1212

# Sheet "全部类型" (all types) - mirrors the same events, update column F too
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 19
$ws4.Range("F6").Value = 20
$ws4.Range("F9").Value = 20
$ws4.Range("F10").Value = 1212
